# chore: update Sheets via scheduled runner
# Refreshes cached market-board pricing/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) on a handful of rows across the crafting-job
# worksheets, as produced by the scheduled data-refresh job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 768.1429000000001
$ws.Range("I33").Value = 811.8461
$ws.Range("K33").Value = 811.8461
$ws.Range("M33").Value = -582.8461
$ws.Range("H62").Value = 3975
$ws.Range("J62").Value = 4300
$ws.Range("L62").Value = 4300
$ws.Range("N62").Value = -5548
$ws.Range("H65").Value = 3975
$ws.Range("J65").Value = 4300
$ws.Range("L65").Value = 21500
$ws.Range("N65").Value = -27740
$ws.Range("H98").Value = 855.95654
$ws.Range("I98").Value = 908.9048
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 908.9048
$ws.Range("L98").Value = 300
$ws.Range("M98").Value = 589.0952
$ws.Range("N98").Value = -3296
$ws.Range("H122").Value = 855.95654
$ws.Range("I122").Value = 908.9048
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 2726.7144
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = -276.7143999999998
$ws.Range("N122").Value = -5800
$ws.Range("H132").Value = 5043.706
$ws.Range("I132").Value = 5889.1562
$ws.Range("J132").Value = 3619.7896
$ws.Range("K132").Value = 17667.4686
$ws.Range("L132").Value = 10859.3688
$ws.Range("M132").Value = -15137.4686
$ws.Range("N132").Value = -15919.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1568.8864
$ws.Range("I61").Value = 1223.0605
$ws.Range("J61").Value = 2606.3635
$ws.Range("K61").Value = 1223.0605
$ws.Range("L61").Value = 2606.3635
$ws.Range("M61").Value = -1011.0605
$ws.Range("N61").Value = -3030.3635
$ws.Range("H132").Value = 2557.2886
$ws.Range("I132").Value = 2371.6904
$ws.Range("J132").Value = 3336.8
$ws.Range("K132").Value = 7115.0712
$ws.Range("L132").Value = 10010.4
$ws.Range("M132").Value = -4585.0712
$ws.Range("N132").Value = -15070.4
$ws.Range("H136").Value = 1568.8864
$ws.Range("I136").Value = 1223.0605
$ws.Range("J136").Value = 2606.3635
$ws.Range("K136").Value = 3669.1815
$ws.Range("L136").Value = 7819.0905
$ws.Range("M136").Value = -1119.1815
$ws.Range("N136").Value = -12919.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1843.6957
$ws.Range("I105").Value = 1456.875
$ws.Range("K105").Value = 1456.875
$ws.Range("M105").Value = 290.125
$ws.Range("H107").Value = 1548.1482
$ws.Range("I107").Value = 1194.2632
$ws.Range("J107").Value = 2388.625
$ws.Range("K107").Value = 1194.2632
$ws.Range("L107").Value = 2388.625
$ws.Range("M107").Value = 725.7367999999999
$ws.Range("N107").Value = -6228.625
$ws.Range("H134").Value = 1139.8182
$ws.Range("I134").Value = 986.8
$ws.Range("J134").Value = 1734.8889
$ws.Range("K134").Value = 2960.4
$ws.Range("L134").Value = 5204.6667
$ws.Range("M134").Value = -425.3999999999996
$ws.Range("N134").Value = -10274.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2768.0425
$ws.Range("I31").Value = 1341.95
$ws.Range("J31").Value = 3824.4075
$ws.Range("K31").Value = 1341.95
$ws.Range("L31").Value = 3824.4075
$ws.Range("M31").Value = -1046.95
$ws.Range("N31").Value = -4414.407499999999
$ws.Range("H34").Value = 2768.0425
$ws.Range("I34").Value = 1341.95
$ws.Range("J34").Value = 3824.4075
$ws.Range("K34").Value = 1341.95
$ws.Range("L34").Value = 3824.4075
$ws.Range("M34").Value = -1139.95
$ws.Range("N34").Value = -4228.407499999999
$ws.Range("H105").Value = 642
$ws.Range("I105").Value = 437.25
$ws.Range("J105").Value = 915
$ws.Range("K105").Value = 437.25
$ws.Range("L105").Value = 915
$ws.Range("M105").Value = 1309.75
$ws.Range("N105").Value = -4409
$ws.Range("H107").Value = 1168.3529
$ws.Range("I107").Value = 983
$ws.Range("J107").Value = 2033.3334
$ws.Range("K107").Value = 983
$ws.Range("L107").Value = 2033.3334
$ws.Range("M107").Value = 937
$ws.Range("N107").Value = -5873.3334
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 1546.2439
$ws.Range("I132").Value = 1075.7931
$ws.Range("J132").Value = 2683.1667
$ws.Range("K132").Value = 3227.379300000001
$ws.Range("L132").Value = 8049.500100000001
$ws.Range("M132").Value = -697.3793000000005
$ws.Range("N132").Value = -13109.5001
$ws.Range("H134").Value = 2857.743
$ws.Range("I134").Value = 3392.8147
$ws.Range("J134").Value = 1051.875
$ws.Range("K134").Value = 10178.4441
$ws.Range("L134").Value = 3155.625
$ws.Range("M134").Value = -7643.444100000001
$ws.Range("N134").Value = -8225.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 973.64703
$ws.Range("I5").Value = 827.94446
$ws.Range("J5").Value = 1137.5625
$ws.Range("K5").Value = 2483.83338
$ws.Range("L5").Value = 3412.6875
$ws.Range("M5").Value = -2371.83338
$ws.Range("N5").Value = -3636.6875
$ws.Range("H17").Value = 600
$ws.Range("J17").Value = 666.6667
$ws.Range("L17").Value = 2000.0001
$ws.Range("N17").Value = -2338.0001
$ws.Range("H107").Value = 155.29167
$ws.Range("I107").Value = 143
$ws.Range("J107").Value = 169.81818
$ws.Range("K107").Value = 429
$ws.Range("L107").Value = 509.4545400000001
$ws.Range("M107").Value = 1491
$ws.Range("N107").Value = -4349.45454
$ws.Range("H131").Value = 924.5517
$ws.Range("I131").Value = 559.8182
$ws.Range("J131").Value = 1147.4445
$ws.Range("K131").Value = 1679.4546
$ws.Range("L131").Value = 3442.3335
$ws.Range("M131").Value = 3360.5454
$ws.Range("N131").Value = -13522.3335
$ws.Range("H132").Value = 1148.8823
$ws.Range("I132").Value = 714.55554
$ws.Range("J132").Value = 1637.5
$ws.Range("K132").Value = 6430.99986
$ws.Range("L132").Value = 14737.5
$ws.Range("M132").Value = -3900.99986
$ws.Range("N132").Value = -19797.5
$ws.Range("H135").Value = 973.64703
$ws.Range("I135").Value = 827.94446
$ws.Range("J135").Value = 1137.5625
$ws.Range("K135").Value = 7451.50014
$ws.Range("L135").Value = 10238.0625
$ws.Range("M135").Value = -4916.50014
$ws.Range("N135").Value = -15308.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1244.4642
$ws.Range("I122").Value = 1281.0869
$ws.Range("K122").Value = 3843.2607
$ws.Range("M122").Value = -1393.2607
$ws.Range("H129").Value = 33333.168
$ws.Range("J129").Value = 33333.168
$ws.Range("L129").Value = 33333.168
$ws.Range("N129").Value = -43333.168
$ws.Range("H132").Value = 3021.8774
$ws.Range("I132").Value = 3072.742
$ws.Range("J132").Value = 2934.2778
$ws.Range("K132").Value = 9218.226000000001
$ws.Range("L132").Value = 8802.8334
$ws.Range("M132").Value = -6688.226000000001
$ws.Range("N132").Value = -13862.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1584.2858
$ws.Range("I40").Value = 1445.3529
$ws.Range("J40").Value = 2174.75
$ws.Range("K40").Value = 1445.3529
$ws.Range("L40").Value = 2174.75
$ws.Range("M40").Value = -1309.3529
$ws.Range("N40").Value = -2446.75
$ws.Range("H122").Value = 3916.7812
$ws.Range("I122").Value = 4679.8887
$ws.Range("J122").Value = 2935.6428
$ws.Range("K122").Value = 14039.6661
$ws.Range("L122").Value = 8806.928400000001
$ws.Range("M122").Value = -11589.6661
$ws.Range("N122").Value = -13706.9284
$ws.Range("H132").Value = 7697454
$ws.Range("I132").Value = 8338525
$ws.Range("J132").Value = 4599.6
$ws.Range("K132").Value = 25015575
$ws.Range("L132").Value = 13798.8
$ws.Range("M132").Value = -25013045
$ws.Range("N132").Value = -18858.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1805.9574
$ws.Range("I132").Value = 709.3
$ws.Range("J132").Value = 8072.5713
$ws.Range("K132").Value = 2127.9
$ws.Range("L132").Value = 24217.7139
$ws.Range("M132").Value = 402.1000000000004
$ws.Range("N132").Value = -29277.7139
